$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update styles in column G to match the rest of the table ---
# Header cell G1 should use the same header style as F1
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)

# Data cells G2:G13 should use the same bordered style as the other data columns (e.g. E)
$ws.Range("E2").Copy()
$ws.Range("G2:G13").PasteSpecial(-4122)

# --- Add the new row (14) of data, cloning formatting from row 13 ---
$ws.Range("A13:G13").Copy()
$ws.Range("A14:G14").PasteSpecial(-4122)

$ws.Range("A14").Value = "Lixol"
$ws.Range("B14").Value = "IPA"
$ws.Range("C14").Value = "Brasil"
$ws.Range("D14").Value = 0.03
$ws.Range("E14").Value = 20
$ws.Range("F14").Value = 0.5
$ws.Range("G14").Value = 23

# --- Update the selected range / active cell to reflect the new extent ---
$null = $ws.Range("G1:G14").Select()
